$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. New shared strings & cell values, in the order needed to match target sst index order:
#    9=sampleBegin 10=column 11=test date text
$ws.Range("B4").Value2 = "sampleBegin"
$ws.Range("C4").Value2 = "column"
$ws.Range("D4").Value2 = "column"
$ws.Range("E4").Value2 = "column"
$ws.Range("E2").Value2 = "test date text"

# 2. Style row 4 header cells: bold, border, center  (-> cellXfs index 3)
$hdr = $ws.Range("B4:E4")
$hdr.Font.Bold = $true
$hdr.Font.Name = "Tahoma"
$hdr.Borders.Weight = 2
$hdr.HorizontalAlignment = -4108

# 3. Style column D data cells (D5:D11): border only (-> cellXfs index 4)
$ws.Range("D5").Borders.Weight = 2
$ws.Range("D6").Borders.Weight = 2
$ws.Range("D7").Borders.Weight = 2
$ws.Range("D8").Borders.Weight = 2
$ws.Range("D9").Borders.Weight = 2
$ws.Range("D10").Borders.Weight = 2
$ws.Range("D11").Borders.Weight = 2

# 4. Style columns B, C, E data cells: border + center (-> cellXfs index 5)
$bce = $ws.Range("B5:B11,C5:C11,E5:E11")
$bce.Borders.Weight = 2
$bce.HorizontalAlignment = -4108

# 5. A2 bold only, no border (-> cellXfs index 6)
$ws.Range("A2").Font.Bold = $true

# 6. Fill remaining empty but styled cells with blank (already covered by border application above, which creates them)
$ws.Range("B10").Value2 = $null
$ws.Range("C10").Value2 = $null
$ws.Range("D10").Value2 = $null
$ws.Range("C11").Value2 = $null
$ws.Range("D11").Value2 = $null

# 7. Column widths
$ws.Columns.Item(2).ColumnWidth = 14.166666666666666
$ws.Columns.Item(4).ColumnWidth = 10.43
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666

# 8. Selection
$ws.Range("E3").Select()
